$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 11 ("Dina" / "Castle Brite" block),
# which pushes the existing rows 11-14 down to rows 14-17 (unchanged content).
$ws.Rows("11:13").Insert()

# Row 11: new "Modesto" / Especial record
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = "2021-12-23"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100103
$ws.Range("H11").Value = "Frutos de hueso (carozo)"
$ws.Range("I11").Value = 100103003
$ws.Range("J11").Value = "Damasco"
$ws.Range("K11").Value = "Modesto"
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 360
$ws.Range("N11").Value = 23000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 23500
$ws.Range("Q11").Value = "$/caja 16 kilos"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 1469
$ws.Range("T11").Value = 16

# Row 12: new "Modesto" / Primera record
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = "2021-12-23"
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100103
$ws.Range("H12").Value = "Frutos de hueso (carozo)"
$ws.Range("I12").Value = 100103003
$ws.Range("J12").Value = "Damasco"
$ws.Range("K12").Value = "Modesto"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("Q12").Value = "$/caja 16 kilos"
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 1344
$ws.Range("T12").Value = 16

# Row 13: new "Modesto" / Segunda record
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = "2021-12-23"
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103003
$ws.Range("J13").Value = "Damasco"
$ws.Range("K13").Value = "Modesto"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 240
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 17500
$ws.Range("Q13").Value = "$/caja 16 kilos"
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 1094
$ws.Range("T13").Value = 16
